$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new "exported" column header
$ws.Range("C3").Value = "exported"

# Move "ingame Card Popup" from B4 to B7
$ws.Range("B4").ClearContents()
$ws.Range("B7").Value = "ingame Card Popup"

# Mark exported (checkmark) cells in column C
$ws.Range("C7").Value = [char]0x2705
$ws.Range("C9").Value = [char]0x2705
$ws.Range("C10").Value = [char]0x2705
$ws.Range("C14").Value = [char]0x2705
$ws.Range("C15").Value = [char]0x2705
$ws.Range("C16").Value = [char]0x2705
$ws.Range("C17").Value = [char]0x2705
$ws.Range("C23").Value = [char]0x2705
$ws.Range("C24").Value = [char]0x2705

# Clear rows 11 and 18 content (SFX_Ingame_Pencil_draw_loop_03 and extra SFX_Ingame_Char_)
$ws.Range("A11").ClearContents()
$ws.Range("A18").ClearContents()

# Update selection to C7
$ws.Range("C7").Select()

$wb.Save()
